$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price), E (Volume(1h)) and G (Hora) store numeric-looking values
# as plain text in this workbook (coinranking.com scrape output). Excel
# auto-converts a numeric-looking string into a real number on assignment, so
# force a text number format on each such cell right before writing it.
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextCell "D2" "296.91"
Set-TextCell "E2" "1.66%"
Set-TextCell "G2" "5"

Set-TextCell "D3" "41.79"
Set-TextCell "E3" "3.64%"
Set-TextCell "G3" "5"

Set-TextCell "D4" "5.006"
Set-TextCell "E4" "-0.34%"
Set-TextCell "G4" "5"

Set-TextCell "D5" "0.07522"
Set-TextCell "E5" "2.74%"
Set-TextCell "G5" "5"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D6" "1.583"
Set-TextCell "E6" "4.23%"
Set-TextCell "G6" "5"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D7" "0.9260"
Set-TextCell "E7" "-0.16%"
Set-TextCell "G7" "5"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D8" "2.401"
Set-TextCell "E8" "0.97%"
Set-TextCell "G8" "5"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D9" "0.1194"
Set-TextCell "E9" "0.27%"
Set-TextCell "G9" "5"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1821"
Set-TextCell "E10" "4.53%"
Set-TextCell "G10" "5"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.08921"
Set-TextCell "E11" "2.65%"
Set-TextCell "G11" "5"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.04081"
Set-TextCell "E12" "-5.65%"
Set-TextCell "G12" "5"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.1048"
Set-TextCell "E13" "-0.61%"
Set-TextCell "G13" "5"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001284"
Set-TextCell "E14" "1.25%"
Set-TextCell "G14" "5"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D15" "0.005929"
Set-TextCell "E15" "-0.13%"
Set-TextCell "G15" "5"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D16" "3.356"
Set-TextCell "E16" "0.53%"
Set-TextCell "G16" "5"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D17" "4.378"
Set-TextCell "E17" "1.87%"
Set-TextCell "G17" "5"

Set-TextCell "E18" "0.74%"
Set-TextCell "G18" "5"

Set-TextCell "D19" "8.143"
Set-TextCell "E19" "2.16%"
Set-TextCell "G19" "5"

Set-TextCell "D20" "0.1391"
Set-TextCell "E20" "0.03%"
Set-TextCell "G20" "5"

Set-TextCell "E21" "10.99%"
Set-TextCell "G21" "5"

Set-TextCell "D22" "0.04097"
Set-TextCell "E22" "3.82%"
Set-TextCell "G22" "5"

Set-TextCell "D23" "0.001266"
Set-TextCell "E23" "0.49%"
Set-TextCell "G23" "5"

Set-TextCell "D24" "0.003902"
Set-TextCell "E24" "3.25%"
Set-TextCell "G24" "5"

Set-TextCell "E25" "-3.93%"
Set-TextCell "G25" "5"

Set-TextCell "G26" "5"

Set-TextCell "G27" "5"

Set-TextCell "G28" "5"

Set-TextCell "G29" "5"

Set-TextCell "G30" "5"

Set-TextCell "G31" "5"

Set-TextCell "G32" "5"

Set-TextCell "G33" "5"

Set-TextCell "G34" "5"

Set-TextCell "G35" "5"

Set-TextCell "G36" "5"

Set-TextCell "G37" "5"

Set-TextCell "D38" "0.02402"
Set-TextCell "E38" "5.07%"
Set-TextCell "G38" "5"

Set-TextCell "D39" "0.05228"
Set-TextCell "E39" "5.02%"
Set-TextCell "G39" "5"

Set-TextCell "D40" "0.006305"
Set-TextCell "E40" "20.03%"
Set-TextCell "G40" "5"

Set-TextCell "D41" "0.007825"
Set-TextCell "E41" "1.45%"
Set-TextCell "G41" "5"

Set-TextCell "E42" "3.31%"
Set-TextCell "G42" "5"

Set-TextCell "D43" "0.007402"
Set-TextCell "E43" "0.73%"
Set-TextCell "G43" "5"

Set-TextCell "D44" "0.007249"
Set-TextCell "E44" "-7.97%"
Set-TextCell "G44" "5"

Set-TextCell "D45" "0.2962"
Set-TextCell "E45" "-6.41%"
Set-TextCell "G45" "5"

Set-TextCell "D46" "0.00006587"
Set-TextCell "E46" "4.22%"
Set-TextCell "G46" "5"

Set-TextCell "D47" "0.00000000751"
Set-TextCell "E47" "-0.03%"
Set-TextCell "G47" "5"

Set-TextCell "D48" "0.04500"
Set-TextCell "E48" "120.43%"
Set-TextCell "G48" "5"

Set-TextCell "D49" "0.004203"
Set-TextCell "E49" "0.04%"
Set-TextCell "G49" "5"

Set-TextCell "D50" "0.00002102"
Set-TextCell "E50" "-0.03%"
Set-TextCell "G50" "5"

Set-TextCell "D51" "0.0002002"
Set-TextCell "E51" "-0.03%"
Set-TextCell "G51" "5"

